$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename existing "Estimate" header and add new production column ---
$ws.Range("B1").Value = "Estimate (review)"
$ws.Range("C1").Value = "Estimate (production)"

# --- Regression coefficient rows (2-7): add production-model estimates in column C ---
$ws.Range("C2").Value = [double]"4.2495693188519699"
$ws.Range("C3").Value = [double]"5.9158296383032896E-3"
$ws.Range("C4").Value = [double]"-1.4879271769975999E-5"
$ws.Range("C5").Value = [double]"4.3176025642281702E-5"
$ws.Range("C6").Value = [double]"-1.11718223359496E-4"
$ws.Range("C7").Value = [double]"4.9118132407661597E-3"

# --- Yearly rows (8-19): add production-time coefficients in column C ---
$ws.Range("C8").Value = [double]"0.101646745352951"
$ws.Range("C9").Value = [double]"-0.12942827953389399"
$ws.Range("C10").Value = [double]"-0.109288525365644"
$ws.Range("C11").Value = [double]"-0.46095013149497499"
$ws.Range("C12").Value = [double]"-0.57019108715497102"
$ws.Range("C13").Value = [double]"-0.56629724829058503"
$ws.Range("C14").Value = [double]"-0.66378718583387797"
$ws.Range("C15").Value = [double]"-0.56771495041387299"
$ws.Range("C16").Value = [double]"-0.47268199440374897"
$ws.Range("C17").Value = [double]"-0.37085600483803999"
$ws.Range("C18").Value = [double]"-0.53778005058935596"
$ws.Range("C19").Value = [double]"-0.24988634171385199"

# --- Apply the 0.00000 number format to the two numeric data columns ---
$ws.Range("B2:C19").NumberFormat = "0.00000"

# --- Column widths to fit the new, wider headers ---
$ws.Range("B:B").ColumnWidth = 16.85546875
$ws.Range("C:C").ColumnWidth = 20.5703125

# --- Selection / used range bookkeeping ---
$ws.Range("A1:C19").Select()
